$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B32 was previously stored as the text "5" - the annotation was corrected
# and is now entered as the real numeric value 5.
$ws.Cells.Item(32, 2).Value = 5

# Append a new annotation row (row 33) for Ying Tang, extending the sheet
# dimension from A1:H32 to A1:H33.
$ws.Cells.Item(33, 1).Value = "Ying Tang"

# politeness_score for this row stays a text value "3" (not re-typed as a
# number like B32 was), so force the cell to Text format before writing it.
$ws.Cells.Item(33, 2).NumberFormat = "@"
$ws.Cells.Item(33, 2).Value = "3"

$ws.Cells.Item(33, 3).Value = "无"
$ws.Cells.Item(33, 4).Value = "DFT"
$ws.Cells.Item(33, 5).Value = "WRI"
$ws.Cells.Item(33, 6).Value = "74483628-1e12-4bb7-acfc-2ccaf38e6d81"
$ws.Cells.Item(33, 7).Value = "HyIFzx-0b_annotated.xlsx"
$ws.Cells.Item(33, 8).Value = "Acronyms are not properly defined."
